$wb = $excel.ActiveWorkbook

# --- Locate existing sheets ---------------------------------------------
$summary = $wb.Worksheets.Item(1)   # "总计"
$q2      = $wb.Worksheets.Item(2)   # "2022-Q2" (current 2nd sheet, will stay 2022-Q2)

# --- Create the new "2022-Q3" sheet by copying the "2022-Q2" sheet ------
# (keeps header row formatting / sheetPr / column layout identical to its
# sibling quarter sheets, just like the original workbook's sheets do)
# Copy placed directly *before* q2 so it lands in the correct final
# position (right after "总计") without needing a separate Move() step.
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# --- Remove the old data rows (2,3,4 inherited from 2022-Q2) except hdr -
$q3.Range("A3:A5").EntireRow.Delete()

# --- Fill in the single fund row for 2022-Q3 ----------------------------
$q3.Cells.Item(2,1).Value = 0

$q3.Cells.Item(2,2).NumberFormat = "@"
$q3.Cells.Item(2,2).Value = "011001"

$q3.Cells.Item(2,3).NumberFormat = "@"
$q3.Cells.Item(2,3).Value = "中邮兴荣价值一年持有期混合"

$q3.Cells.Item(2,4).NumberFormat = "@"
$q3.Cells.Item(2,4).Value = "5.15"

$q3.Cells.Item(2,5).NumberFormat = "@"
$q3.Cells.Item(2,5).Value = "40.76"

$q3.Cells.Item(2,6).NumberFormat = "@"
$q3.Cells.Item(2,6).Value = "2.24"

$q3.Cells.Item(2,7).NumberFormat = "@"
$q3.Cells.Item(2,7).Value = "0.1154"

$q3.Cells.Item(2,8).Value = 9

# --- Update the "总计" summary sheet: insert a new row for 2022-Q3 ------
$summary.Rows.Item(2).Insert()
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)   # xlPasteFormats

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 1
$summary.Cells.Item(2,4).Value = 0.12
